$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Estados Unidos (row 4) stats ---
$ws.Range("B4").Value = 1179454
$ws.Range("C4").Value = 18680
$ws.Range("D4").Value = 166203
$ws.Range("E4").Value = 945072
$ws.Range("G4").Value = 735
$ws.Range("H4").Value = 68179

# --- Update Irlanda (row 26) new-deaths figure ---
$ws.Range("F26").Value = 98

# --- Costa Rica moves up (new case total overtakes Libano & Niger) ---
# Insert a fresh row right before "Libano" (row 99) and populate it with
# Costa Rica's updated figures; then remove Costa Rica's old row (which,
# after the insert shifts everything down by one, now sits at row 102).
$ws.Rows(99).Insert()
$ws.Range("A99").Value = "Costa Rica"
$ws.Range("B99").Value = 739
$ws.Range("C99").Value = 6
$ws.Range("D99").Value = 386
$ws.Range("E99").Value = 347
$ws.Range("F99").Value = 6
$ws.Range("G99").Value = 0
$ws.Range("H99").Value = 6
$ws.Rows(102).Delete()
